# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the file that was
# "85aead02-2c32-45ba-ad3f-dcd9c4f70651.md" (previously "In Translation")
# is now "Ready for handoff", while "b9dfe6b1-4627-49bb-b78d-db46cadef183.md"
# stays "In Translation". The two rows in every sheet swap positions (the
# b9dfe6b1 row moves up to row 2, the 85aead row moves down to row 3), and
# the 85aead row picks up fresh handoff metadata (new datetime, priority
# "mt" instead of "ht").

$wb = $excel.ActiveWorkbook

$url85aead   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74d625e46884e6804dbe078972837b90c5448650/e2e/85aead02-2c32-45ba-ad3f-dcd9c4f70651.md"
$urlB9dfe6b1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74d625e46884e6804dbe078972837b90c5448650/e2e/b9dfe6b1-4627-49bb-b78d-db46cadef183.md"

$name85aead   = "85aead02-2c32-45ba-ad3f-dcd9c4f70651.md"
$nameB9dfe6b1 = "b9dfe6b1-4627-49bb-b78d-db46cadef183.md"

# ----------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 2 now carries the b9dfe6b1 file, row 3 the 85aead file (swapped).
$ws1.Range("A2").Value = $nameB9dfe6b1
$ws1.Range("A3").Value = $name85aead

# Status / datetime for row 3 (85aead) move to "Ready for handoff".
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-09-05 20:17:59"

# Rebuild the two hyperlinks: targets (r:id -> URL) stay exactly as they
# were, only the display text (and the cell it decorates) swaps.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $url85aead,   "", "", "e2e\" + $nameB9dfe6b1)
$ws1.Hyperlinks.Add($ws1.Range("B3"), $urlB9dfe6b1, "", "", "e2e\" + $name85aead)

# Columns E/F widen to fit "Ready for handoff".
$ws1.Columns.Item(5).ColumnWidth = 17.2159881591797
$ws1.Columns.Item(6).ColumnWidth = 17.2159881591797

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $nameB9dfe6b1
$ws2.Range("G2").Value = "b9dfe6b1-4627-49bb-b78d-db46cadef183.3dff16e31c78f42ac0b0927139d2c4b1b73b996d.zh-cn.xlf"

$ws2.Range("A3").Value = $name85aead
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("E3").Value = "mt"
$ws2.Range("G3").Value = "85aead02-2c32-45ba-ad3f-dcd9c4f70651.e1717f842d4110fcd965b06e615eb01f8aae51dc.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-09-05 20:17:55"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $url85aead,   "", "", $nameB9dfe6b1)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $urlB9dfe6b1, "", "", $name85aead)

$ws2.Columns.Item(3).ColumnWidth = 17.2159881591797

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $nameB9dfe6b1
$ws3.Range("G2").Value = "b9dfe6b1-4627-49bb-b78d-db46cadef183.3dff16e31c78f42ac0b0927139d2c4b1b73b996d.de-de.xlf"

$ws3.Range("A3").Value = $name85aead
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("E3").Value = "mt"
$ws3.Range("G3").Value = "85aead02-2c32-45ba-ad3f-dcd9c4f70651.e1717f842d4110fcd965b06e615eb01f8aae51dc.de-de.xlf"
$ws3.Range("H3").Value = "2016-09-05 20:17:59"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $url85aead,   "", "", $nameB9dfe6b1)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $urlB9dfe6b1, "", "", $name85aead)

$ws3.Columns.Item(3).ColumnWidth = 17.2159881591797
